# Remove the "Implement AutoReplay feature, to rollback generated mock
# after creating" bullet paragraph entirely (including its paragraph
# mark), since the feature described by it has been implemented and the
# to-do item is no longer needed.

$d = $word.ActiveDocument

# Locate the paragraph whose text contains "AutoReplay" and delete its
# whole range (text + the trailing paragraph mark), which collapses it
# out of the list rather than leaving an empty bullet behind.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*AutoReplay*") {
        $para.Range.Delete()
    }
}
